$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 130, shifting existing rows 130-212 down to 132-214
$ws.Rows("130:131").Insert()

# New row 130 data
$ws.Cells.Item(130, 1).Value = 9
$ws.Cells.Item(130, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(130, 3).Value = "Metropolitana"
$ws.Cells.Item(130, 4).Value = 44574
$ws.Cells.Item(130, 5).Value = 13
$ws.Cells.Item(130, 6).Value = 100112043
$ws.Cells.Item(130, 7).Value = "Pepino ensalada"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 79
$ws.Cells.Item(130, 11).Value = 9000
$ws.Cells.Item(130, 12).Value = 10000
$ws.Cells.Item(130, 13).Value = 9494
$ws.Cells.Item(130, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(130, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(130, 16).Value = 158
$ws.Cells.Item(130, 17).Value = 60
$ws.Cells.Item(130, 18).Value = "Hortaliza"

# New row 131 data
$ws.Cells.Item(131, 1).Value = 9
$ws.Cells.Item(131, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(131, 3).Value = "Metropolitana"
$ws.Cells.Item(131, 4).Value = 44574
$ws.Cells.Item(131, 5).Value = 13
$ws.Cells.Item(131, 6).Value = 100112043
$ws.Cells.Item(131, 7).Value = "Pepino ensalada"
$ws.Cells.Item(131, 8).Value = "Sin especificar"
$ws.Cells.Item(131, 9).Value = "Primera"
$ws.Cells.Item(131, 10).Value = 106
$ws.Cells.Item(131, 11).Value = 11000
$ws.Cells.Item(131, 12).Value = 12000
$ws.Cells.Item(131, 13).Value = 11500
$ws.Cells.Item(131, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(131, 15).Value = "Región del Maule"
$ws.Cells.Item(131, 16).Value = 192
$ws.Cells.Item(131, 17).Value = 60
$ws.Cells.Item(131, 18).Value = "Hortaliza"
